$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 37.44140625
$ws.Columns.Item(3).ColumnWidth = 10.33203125

# Row 1 - title row
$ws.Range("A1").Value = "Energieverbruik CO2 vervloeiing"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Underline = $true

# Row 2 - headers
$ws.Range("A2").Value = "Energie verbruik totaal maand"
$ws.Range("B2").Value = "##"
$ws.Range("C2").Value = "MWh"

# Row 3 - data row
$ws.Range("A3").Value = "Energieverbruik (garantie 0,14 kWh/Nm3 biogas)"
$ws.Range("B3").Value = "##"
$ws.Range("C3").Value = "kWh/Nm3 biogas"
$ws.Range("A3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("A3").RowHeight = 28.8

# Gray fill on B2 and B3
$ws.Range("B2").Interior.ThemeColor = 1
$ws.Range("B2").Interior.TintAndShade = -0.14999847407452621
$ws.Range("B3").Interior.ThemeColor = 1
$ws.Range("B3").Interior.TintAndShade = -0.14999847407452621

# Borders - edges only, per cell
$ws.Range("A1").Borders.Item(7).LineStyle = 1
$ws.Range("A1").Borders.Item(7).Weight = 2
$ws.Range("A1").Borders.Item(8).LineStyle = 1
$ws.Range("A1").Borders.Item(8).Weight = 2

$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("B1").Borders.Item(8).Weight = 2

$ws.Range("C1").Borders.Item(8).LineStyle = 1
$ws.Range("C1").Borders.Item(8).Weight = 2
$ws.Range("C1").Borders.Item(10).LineStyle = 1
$ws.Range("C1").Borders.Item(10).Weight = 2

$ws.Range("A2").Borders.Item(7).LineStyle = 1
$ws.Range("A2").Borders.Item(7).Weight = 2

$ws.Range("C2").Borders.Item(10).LineStyle = 1
$ws.Range("C2").Borders.Item(10).Weight = 2

$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(7).Weight = 2
$ws.Range("A3").Borders.Item(9).LineStyle = 1
$ws.Range("A3").Borders.Item(9).Weight = 2

$ws.Range("B3").Borders.Item(9).LineStyle = 1
$ws.Range("B3").Borders.Item(9).Weight = 2

$ws.Range("C3").Borders.Item(9).LineStyle = 1
$ws.Range("C3").Borders.Item(9).Weight = 2
$ws.Range("C3").Borders.Item(10).LineStyle = 1
$ws.Range("C3").Borders.Item(10).Weight = 2

$ws.Range("A1:C3").Select()
